$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'246.29"
$ws.Range("D4").Value2 = "'5.315"
$ws.Range("D5").Value2 = "'0.05883"
$ws.Range("D6").Value2 = "'3.391"
$ws.Range("D8").Value2 = "'0.8136"
$ws.Range("D9").Value2 = "'0.9568"
$ws.Range("B10").Value2 = "WazirX"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value2 = "'0.1417"
$ws.Range("E10").Value2 = "9WazirXWRX"
$ws.Range("B11").Value2 = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value2 = "'0.03735"
$ws.Range("E11").Value2 = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value2 = "MandalaExchangeToken"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value2 = "'0.07331"
$ws.Range("E12").Value2 = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value2 = "BitrueCoin"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value2 = "'0.03054"
$ws.Range("E13").Value2 = "12BitrueCoinBTR"
$ws.Range("B14").Value2 = "MCDex"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value2 = "'4.414"
$ws.Range("E14").Value2 = "13MCDexMCB"
$ws.Range("B15").Value2 = "BitMartToken"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value2 = "'0.09395"
$ws.Range("E15").Value2 = "14BitMartTokenBMX"
$ws.Range("B16").Value2 = "BitForexToken"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value2 = "'0.001594"
$ws.Range("E16").Value2 = "15BitForexTokenBF"
$ws.Range("B17").Value2 = "CoinExToken"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value2 = "'0.04810"
$ws.Range("E17").Value2 = "16CoinExTokenCET"
$ws.Range("B18").Value2 = "One"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value2 = "'0.0005902"
$ws.Range("E18").Value2 = "17OneONE"
$ws.Range("D19").Value2 = "'0.006125"
$ws.Range("D20").Value2 = "'0.004085"
$ws.Range("D21").Value2 = "'0.0009839"
$ws.Range("D23").Value2 = "'3.684"
$ws.Range("D41").Value2 = "'0.006612"
$ws.Range("D42").Value2 = "'0.1074"
$ws.Range("D44").Value2 = "'0.005907"
$ws.Range("D45").Value2 = "'0.00005670"
$ws.Range("D48").Value2 = "'0.002671"
